$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rewrite C1 and insert new headers, reusing the
# bold/centered/bordered style already applied to the existing header cells.
$ws.Range("A1").Copy()
$ws.Range("J1:L1").PasteSpecial(-4122)

$ws.Range("C1").Value = "original error rate PC_LabelCorrection"
$ws.Range("D1").Value = "error rate after correction PC_LabelCorrection"
$ws.Range("E1").Value = "original error rate CL"
$ws.Range("F1").Value = "error rate after correction CL"
$ws.Range("G1").Value = "before_fix_mean"
$ws.Range("H1").Value = "before_fix_variance"
$ws.Range("I1").Value = "before_fix_std"
$ws.Range("J1").Value = "after_fix_mean"
$ws.Range("K1").Value = "after_fix_variance"
$ws.Range("L1").Value = "after_fix_std"

# --- Data rows: C:F now hold the error-rate numbers, G:I hold the
# before_fix_* stats and J:L hold the after_fix_* stats. Rows that don't
# carry a value for a block are cleared back to blank.

function Set-Row($r, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l) {
    if ($c -eq $null) { $ws.Range("C$r").ClearContents() } else { $ws.Range("C$r").Value = $c }
    if ($d -eq $null) { $ws.Range("D$r").ClearContents() } else { $ws.Range("D$r").Value = $d }
    if ($e -eq $null) { $ws.Range("E$r").ClearContents() } else { $ws.Range("E$r").Value = $e }
    if ($f -eq $null) { $ws.Range("F$r").ClearContents() } else { $ws.Range("F$r").Value = $f }
    if ($g -eq $null) { $ws.Range("G$r").ClearContents() } else { $ws.Range("G$r").Value = $g }
    if ($h -eq $null) { $ws.Range("H$r").ClearContents() } else { $ws.Range("H$r").Value = $h }
    if ($i -eq $null) { $ws.Range("I$r").ClearContents() } else { $ws.Range("I$r").Value = $i }
    if ($j -eq $null) { $ws.Range("J$r").ClearContents() } else { $ws.Range("J$r").Value = $j }
    if ($k -eq $null) { $ws.Range("K$r").ClearContents() } else { $ws.Range("K$r").Value = $k }
    if ($l -eq $null) { $ws.Range("L$r").ClearContents() } else { $ws.Range("L$r").Value = $l }
}

# row 2: metric_2D / PC
Set-Row 2 0.08169 0.04607 0.099 0.09702 $null $null $null $null $null $null
# row 3: metric_breast_cancer / PC
Set-Row 3 0.08225 0.03779 0.0984 0.04482 $null $null $null $null $null $null
# row 4: metric_load_iris / PC
Set-Row 4 0.05734 0.01866 0.1 0.058 $null $null $null $null $null $null
# row 5: metric_load_wine / PC
Set-Row 5 0.07417 0.02248 0.0899 0.03371 $null $null $null $null $null $null
# row 6: indices_PC_LabelCorrection_before_fix_OCPC / PC
Set-Row 6 $null $null $null $null 0.0738625 0.00010117686875 0.01005867132130283 $null $null $null
# row 7: indices_CL_before_fix_OCPC / PC
Set-Row 7 $null $null $null $null 0.09682499999999999 0.00001631187500000003 0.004038796231552173 $null $null $null
# row 8: indices_PC_LabelCorrection_after_fix_OCPC / PC
Set-Row 8 $null $null $null $null $null $null $null 0.02083333333333333 0.0002999847222222222 0.01732006703861801
# row 9: indices_CL_after_fix_OCPC / PC
Set-Row 9 $null $null $null $null $null $null $null 0.03892499999999999 0.001138523858333333 0.03374201917984952
# row 10: metric_2D / LOF
Set-Row 10 0.06832000000000001 0.05545 0.099 0.09801 $null $null $null $null $null $null
# row 11: metric_breast_cancer / LOF
Set-Row 11 0.09068999999999999 0.0355 0.0984 0.04569 $null $null $null $null $null $null
# row 12: metric_load_iris / LOF
Set-Row 12 0.11133 0.03599 0.1 0.06399000000000001 $null $null $null $null $null $null
# row 13: metric_load_wine / LOF
Set-Row 13 0.12304 0.03706 0.0899 0.03482 $null $null $null $null $null $null
# row 14: indices_PC_LabelCorrection_before_fix_LOF / LOF
Set-Row 14 $null $null $null $null 0.09834499999999999 0.0004346382249999999 0.02084797891883048 $null $null $null
# row 15: indices_CL_before_fix_LOF / LOF
Set-Row 15 $null $null $null $null 0.09682499999999999 0.00001631187500000003 0.004038796231552173 $null $null $null
# row 16: indices_PC_LabelCorrection_after_fix_LOF / LOF
Set-Row 16 $null $null $null $null $null $null $null 0.02733333333333333 0.0004201682555555556 0.02049800613609908
# row 17: indices_CL_after_fix_LOF / LOF
Set-Row 17 $null $null $null $null $null $null $null 0.04041833333333333 0.001199806447222222 0.03463822234500816
